$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted right after row 327 ("Fruta / hortaliza, semanal").
# This pushes every existing record from row 328 through row 456 down by one row
# (to rows 329 through 457), and row 328 is populated with the new data point.
$ws.Rows.Item(328).Insert()

$ws.Cells.Item(328, 1).Value = 5
$ws.Cells.Item(328, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(328, 3).Value = "Maule"
$ws.Cells.Item(328, 4).Value = 45009
$ws.Cells.Item(328, 5).Value = 7
$ws.Cells.Item(328, 6).Value = 100112003
$ws.Cells.Item(328, 7).Value = "Ajo"
$ws.Cells.Item(328, 8).Value = "Chino"
$ws.Cells.Item(328, 9).Value = "Primera"
$ws.Cells.Item(328, 10).Value = 200
$ws.Cells.Item(328, 11).Value = 19000
$ws.Cells.Item(328, 12).Value = 19000
$ws.Cells.Item(328, 13).Value = 19000
$ws.Cells.Item(328, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(328, 15).Value = "China"
$ws.Cells.Item(328, 16).Value = 1900
$ws.Cells.Item(328, 17).Value = 10
$ws.Cells.Item(328, 18).Value = "Hortaliza"
